# "Le cost c'est fun" - fix the assembly header IDs in the WT BOM sheet so
# they match the actual assembly numbers (WT_A0100 / WT_A0200 / WT_A0300)
# instead of the generic placeholders (WT_A0001 / WT_A0002 / WT_A0003), and
# restore the view to the top of the sheet with G6 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G17").Value = "WT_A0300"
$ws.Range("G2").Value = "WT_A0100"
$ws.Range("G5").Value = "WT_A0200"

$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("G6").Select()
